$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.292394999999999
$ws.Range("H2").Value = 21.877185
$ws.Range("I2").Value = 0.5244715940033005
$ws.Range("J2").Value = 0.5244715940033005
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.391302333333333
$ws.Range("N2").Value = 4.173907
$ws.Range("O2").Value = 0.03171126955348368
$ws.Range("P2").Value = 0.03171126955348368
$ws.Range("Q2").Value = 10.14592617908833
$ws.Range("R2").Value = 91.31333561179498
$ws.Range("S2").Value = 0.01663166009058391
$ws.Range("T2").Value = 0.01663166009058391
$ws.Range("G3").Value = 7.292394999999999
$ws.Range("H3").Value = 21.877185
$ws.Range("I3").Value = 0.5244715940033005
$ws.Range("J3").Value = 0.5244715940033005
$ws.Range("O3").Value = 0.9045117298527411
$ws.Range("P3").Value = 0.9045117298527411
$ws.Range("Q3").Value = 289.3958320945633
$ws.Range("R3").Value = 2604.56248885107
$ws.Range("S3").Value = 0.4743907087505498
$ws.Range("T3").Value = 0.4743907087505498
$ws.Range("G4").Value = 7.292394999999999
$ws.Range("H4").Value = 21.877185
$ws.Range("I4").Value = 0.5244715940033005
$ws.Range("J4").Value = 0.5244715940033005
$ws.Range("O4").Value = 0.06377700059377522
$ws.Range("P4").Value = 0.06377700059377524
$ws.Range("Q4").Value = 20.40526125441833
$ws.Range("R4").Value = 183.647351289765
$ws.Range("S4").Value = 0.03344922516216673
$ws.Range("T4").Value = 0.03344922516216674
$ws.Range("I5").Value = 0.02354183170388992
$ws.Range("J5").Value = 0.02354183170388992
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.391302333333333
$ws.Range("N5").Value = 4.173907
$ws.Range("O5").Value = 0.03171126955348368
$ws.Range("P5").Value = 0.03171126955348368
$ws.Range("Q5").Value = 0.4554177753746667
$ws.Range("R5").Value = 4.098759978372001
$ws.Range("S5").Value = 0.0007465413709448013
$ws.Range("T5").Value = 0.0007465413709448013
$ws.Range("I6").Value = 0.02354183170388992
$ws.Range("J6").Value = 0.02354183170388992
$ws.Range("O6").Value = 0.9045117298527411
$ws.Range("P6").Value = 0.9045117298527411
$ws.Range("S6").Value = 0.02129386291838758
$ws.Range("T6").Value = 0.02129386291838758
$ws.Range("I7").Value = 0.02354183170388992
$ws.Range("J7").Value = 0.02354183170388992
$ws.Range("O7").Value = 0.06377700059377522
$ws.Range("P7").Value = 0.06377700059377524
$ws.Range("S7").Value = 0.001501427414557544
$ws.Range("T7").Value = 0.001501427414557544
$ws.Range("I8").Value = 0.4519865742928097
$ws.Range("J8").Value = 0.4519865742928096
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.391302333333333
$ws.Range("N8").Value = 4.173907
$ws.Range("O8").Value = 0.03171126955348368
$ws.Range("P8").Value = 0.03171126955348368
$ws.Range("Q8").Value = 8.743700267368554
$ws.Range("R8").Value = 78.693302406317
$ws.Range("S8").Value = 0.01433306809195496
$ws.Range("T8").Value = 0.01433306809195496
$ws.Range("I9").Value = 0.4519865742928097
$ws.Range("J9").Value = 0.4519865742928096
$ws.Range("O9").Value = 0.9045117298527411
$ws.Range("P9").Value = 0.9045117298527411
$ws.Range("S9").Value = 0.4088271581838038
$ws.Range("T9").Value = 0.4088271581838037
$ws.Range("I10").Value = 0.4519865742928097
$ws.Range("J10").Value = 0.4519865742928096
$ws.Range("O10").Value = 0.06377700059377522
$ws.Range("P10").Value = 0.06377700059377524
$ws.Range("S10").Value = 0.02882634801705095
$ws.Range("T10").Value = 0.02882634801705095
